# Populate the TECHNICAL DETAILS table (Table 1) - rows 1..4, second column.
$d = $word.ActiveDocument

$tech = $d.Tables.Item(1)

$standardPrepText = "It is recommended that the standards be prepared no more than 2 hours prior to performing the experiment. Use one 10 ng of lyophilized Mouse Klk1 standard for each experiment. Gently spin the vial prior to use. Reconstitute the standard to a stock concentration of 10 ng/ml using 1ml of sample diluent. Allow the standard to sit for a minimum of 10 minutes with gentle agitation prior to making dilutions."

# Row 1: Capture/Detection Antibodies
$tech.Cell(1, 2).Range.Text = $standardPrepText

# Row 2: Specificity
$tech.Cell(2, 2).Range.Text = "Natural and recombinant Mouse Klk1"

# Row 3: Standard Protein
$tech.Cell(3, 2).Range.Text = $standardPrepText

# Row 4: Cross-reactivity
$tech.Cell(4, 2).Range.Text = "This kit is for the detection of Mouse Klk1. No significant cross-reactivity or interference between Klk1 and its analogs was observed. This claim is limited by existing techniques; therefore, cross- reactivity may exist with untested analogs."

# --- Pass 1: write the correct final values into the correct cells, scoped per-cell so
#     duplicate text elsewhere in the document cannot cause mistargeting. ---

# Intra-assay precision table (Table 5): Sample / n / Mean(pg/ml) / Standard Deviation / CV(%)
$intra = $d.Tables.Item(5)
$intra.Cell(2, 4).Range.Text = "10.15"   # Sample 1 (unchanged)
$intra.Cell(2, 5).Range.Text = "7.0%"    # Sample 1 (unchanged)
$intra.Cell(3, 4).Range.Text = "11.23"   # Sample 2: 23.03 -> 11.23
$intra.Cell(3, 5).Range.Text = "7.5%"    # Sample 2: 7.0% -> 7.5%
$intra.Cell(4, 4).Range.Text = "9.88"    # Sample 3: 65.84 -> 9.88
$intra.Cell(4, 5).Range.Text = "6.7%"    # Sample 3: 6.2% -> 6.7%

# Inter-assay precision table (Table 6): Sample / n / Mean(pg/ml) / Standard Deviation / CV(%)
$inter = $d.Tables.Item(6)
$inter.Cell(2, 4).Range.Text = "13.05"   # Sample 1 (unchanged)
$inter.Cell(2, 5).Range.Text = "9.0%"    # Sample 1 (unchanged)
$inter.Cell(3, 4).Range.Text = "14.27"   # Sample 2: 29.61 -> 14.27
$inter.Cell(3, 5).Range.Text = "9.6%"    # Sample 2: 9.0% -> 9.6%
$inter.Cell(4, 4).Range.Text = "12.69"   # Sample 3: 95.58 -> 12.69
$inter.Cell(4, 5).Range.Text = "8.8%"    # Sample 3: 9.0% -> 8.8%

# --- Pass 2: each value written above is now unique across the whole document, so a
#     plain document-wide Find/Replace of a value with itself rebuilds a clean run
#     (drops the stale xml:space="preserve" some of these runs originally carried). ---
$cleanupValues = @("10.15", "7.0%", "11.23", "7.5%", "9.88", "6.7%", "13.05", "9.0%", "14.27", "9.6%", "12.69", "8.8%")
foreach ($val in $cleanupValues) {
    $d.Content.Find.Execute($val, $false, $false, $false, $false, $false, $true, 1, $false, $val, 2)
}

Write-Output "Edits applied."
